$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new component row (row 4)
$ws.Range("D4").Value = "PCA9306D"
$ws.Range("B4").Value = "5 to 3,3 bidir i²c"
$ws.Range("C4").Value = 1
$ws.Range("F4").Value = 0.764

# Add the hyperlink on E4 (Excel auto-creates the "Lien hypertexte" style/font)
$url = "http://fr.farnell.com/nxp/pca9306d/ic-i2c-smbus-volt-trans-double/dp/2212070?ost=2212070&mckv=sS7zBkBKq_dc%7Cpcrid%7C79324297994%7Ckword%7Cpca9306d%7Cmatch%7Cp%7Cplid%7C&CMP=KNC-GFR-FFR-GEN-SKU-MDC&gclid=CIzq7vvXm8gCFRITGwodfYMC8w"
$ws.Hyperlinks.Add($ws.Range("E4"), $url)

# Update the active cell selection to F6 (matches the recorded edit)
$ws.Range("F6").Select()
